$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to be treated as text before writing the new
# values, since many of the values look like numbers (e.g. "0.9962") and
# Excel would otherwise silently convert them to floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# --- Rows whose Coin/Link stay the same; only Price (D) and Volume(1h) (E) change ---

$ws.Range("D2").Value = "26.953.71"
$ws.Range("E2").Value = "  +2.14%  "

$ws.Range("D3").Value = "1.736.83"
$ws.Range("E3").Value = "  +0.79%  "

$ws.Range("D4").Value = "0.9962"
$ws.Range("E4").Value = "  -0.33%  "

$ws.Range("D5").Value = "243.46"
$ws.Range("E5").Value = "  +0.29%  "

$ws.Range("D6").Value = "0.9966"
$ws.Range("E6").Value = "  -0.31%  "

$ws.Range("D7").Value = "0.4934"
$ws.Range("E7").Value = "  +1.79%  "

$ws.Range("D8").Value = "0.2627"
$ws.Range("E8").Value = "  +0.98%  "

$ws.Range("D9").Value = "0.06238"
$ws.Range("E9").Value = "  +0.67%  "

$ws.Range("D10").Value = "1.740.61"
$ws.Range("E10").Value = "  +1.01%  "

$ws.Range("D11").Value = "16.12"
$ws.Range("E11").Value = "  +3.90%  "

$ws.Range("D12").Value = "0.06947"
$ws.Range("E12").Value = "  -0.50%  "

$ws.Range("D13").Value = "0.6157"
$ws.Range("E13").Value = "  +3.19%  "

$ws.Range("D14").Value = "4.537"
$ws.Range("E14").Value = "  -0.18%  "

$ws.Range("D15").Value = "77.87"
$ws.Range("E15").Value = "  +0.91%  "

$ws.Range("D16").Value = "0.9963"
$ws.Range("E16").Value = "  -0.34%  "

$ws.Range("D17").Value = "26.683.05"
$ws.Range("E17").Value = "  +1.06%  "

$ws.Range("D18").Value = "0.9962"
$ws.Range("E18").Value = "  -0.36%  "

$ws.Range("D19").Value = "0.000007191"
$ws.Range("E19").Value = "  +0.10%  "

$ws.Range("D20").Value = "11.52"
$ws.Range("E20").Value = "  +1.88%  "

$ws.Range("D21").Value = "1.961.31"
$ws.Range("E21").Value = "  +0.95%  "

$ws.Range("D22").Value = "4.475"
$ws.Range("E22").Value = "  +0.36%  "

$ws.Range("D23").Value = "8.624"
$ws.Range("E23").Value = "  +0.99%  "

$ws.Range("D24").Value = "5.147"
$ws.Range("E24").Value = "  -0.06%  "

$ws.Range("D25").Value = "138.49"
$ws.Range("E25").Value = "  +1.12%  "

$ws.Range("D26").Value = "15.43"
$ws.Range("E26").Value = "  +1.44%  "

$ws.Range("D27").Value = "1.794"
$ws.Range("E27").Value = "  +4.86%  "

# --- Row 28 & 29 swap places (Toncoin now ranks above BitcoinCash) ---

$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "1.393"
$ws.Range("E28").Value = "  -0.68%  "

$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").Value = "106.90"
$ws.Range("E29").Value = "  -0.04%  "

# --- Rows 30-49: only Price (D) and Volume(1h) (E) change ---

$ws.Range("D30").Value = "3.971"
$ws.Range("E30").Value = "  +0.11%  "

$ws.Range("D31").Value = "0.08009"
$ws.Range("E31").Value = "  +0.73%  "

$ws.Range("D32").Value = "3.698"
$ws.Range("E32").Value = "  +0.65%  "

$ws.Range("D33").Value = "0.04526"
$ws.Range("E33").Value = "  +0.47%  "

$ws.Range("D34").Value = "2.618"
$ws.Range("E34").Value = "  +0.65%  "

$ws.Range("D35").Value = "1.016"
$ws.Range("E35").Value = "  +1.89%  "

$ws.Range("D36").Value = "0.6281"
$ws.Range("E36").Value = "  +0.78%  "

$ws.Range("D37").Value = "0.9488"
$ws.Range("E37").Value = "  +5.08%  "

$ws.Range("D38").Value = "2.072"
$ws.Range("E38").Value = "  +6.13%  "

$ws.Range("D39").Value = "2.430"
$ws.Range("E39").Value = "  +1.58%  "

$ws.Range("D40").Value = "0.9955"
$ws.Range("E40").Value = "  -0.43%  "

$ws.Range("D41").Value = "0.01513"
$ws.Range("E41").Value = "  +2.54%  "

$ws.Range("D42").Value = "5.623"
$ws.Range("E42").Value = "  +3.08%  "

$ws.Range("D43").Value = "99.61"
$ws.Range("E43").Value = "  +0.39%  "

$ws.Range("D44").Value = "0.3891"
$ws.Range("E44").Value = "  +1.66%  "

$ws.Range("D45").Value = "7.012"
$ws.Range("E45").Value = "  +4.45%  "

$ws.Range("D46").Value = "0.1168"
$ws.Range("E46").Value = "  +2.02%  "

$ws.Range("D47").Value = "0.05393"
$ws.Range("E47").Value = "  +0.68%  "

$ws.Range("D48").Value = "7.969"
$ws.Range("E48").Value = "  +3.98%  "

$ws.Range("D49").Value = "30.33"
$ws.Range("E49").Value = "  +0.94%  "

# --- Row 50 & 51 swap places (NEARProtocol now ranks above Aave) ---

$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "1.247"
$ws.Range("E50").Value = "  +0.59%  "

$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "52.03"
$ws.Range("E51").Value = "  +2.20%  "

# Restore the Price column's cell style to the workbook default now that the
# text values have been written (undoes the temporary "@" number format).
$ws.Range("D2:D51").Style = "Normal"
